$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: date moved, description added (chapter 10 started) ---
$ws.Range("A9").Value = "10/27/2025"
$ws.Range("E9").Value = "Started chapter 10 which is about modifying files on a computer. This is something that may actually be useful. "
$ws.Rows.Item(9).RowHeight = 28.8

# --- Row 10: date moved, description replaced (chapter 10 finished, longer writeup) ---
$ws.Range("A10").Value = "10/27/2025"
$ws.Range("E10").Value = "Finished chapter 10. One of the practice programs wanted me to use Regex again. I did this in a previous chapter, however I completely forgot how to read that format. Apparently Regex is famous for this. I therefore employed the help of both google search and AI tips to help me in this. This allowed me to work through that and finish the program. Before that I also found another, very inefficient, way of doing the search. I was proud of figuring this out myself, but also realised that it would not be practical."
$ws.Rows.Item(10).RowHeight = 86.4

# --- Row 11: date moved, description now "Almost finished chapter 11" (default row height) ---
$ws.Range("A11").Value = "11/3/2025"
$ws.Range("E11").Value = "Almost finished chapter 11"

# --- Row 12: date moved, description replaced (chapter 11 finished, longer writeup) ---
$ws.Range("A12").Value = "11/4/2025"
$ws.Range("E12").Value = "Finished chapter 11. This chapter was the continuation of the previous chapter. I feel like I can now do some actual automation. From here I can freely choose the next chapter, depending on what I find interesting. I will probably start by doing the excel Chapter."
$ws.Rows.Item(12).RowHeight = 43.2

# --- Shrink Table1 so it no longer covers the trailing blank rows 13-14 ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E12"))

# Rows 13-14 fall outside the table now, so their calculated Time column
# formula goes away, leaving plain (empty) formatted cells.
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()

# --- Column A got a touch wider (best fit) ---
$ws.Columns.Item(1).ColumnWidth = 10.33

# --- View state tweaks (zoom + scroll position + active selection) ---
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("O6").Select() | Out-Null
